$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("L_L")

# --- Clear the existing AutoFilter criterion (SPECIES = YFT) and unhide the rows it hid ---
# This clears sheetPr filterMode / row-level hidden flags while keeping the
# autoFilter range itself (A1:H27) intact.
$ws.ShowAllData()

# --- Append the 2 new length-length parameter rows for ALB (row 28 & 29) ---
# Row 28: ALB / LDF / FL
$ws.Cells.Item(28, 1).Value = "ALB"
$ws.Cells.Item(28, 2).Value = "LDF"
$ws.Cells.Item(28, 3).Value = "FL"
$ws.Cells.Item(28, 4).Value = 1
$ws.Cells.Item(28, 4).NumberFormat = $ws.Cells.Item(2, 4).NumberFormat
$ws.Cells.Item(28, 5).Value = 0.2678
$ws.Cells.Item(28, 6).Value = 5.4938
$ws.Cells.Item(28, 7).Value = "INVPROP"
$ws.Cells.Item(28, 8).Value = "Dhurmeea2016"

# Row 29: ALB / PCL / FL
$ws.Cells.Item(29, 1).Value = "ALB"
$ws.Cells.Item(29, 2).Value = "PCL"
$ws.Cells.Item(29, 3).Value = "FL"
$ws.Cells.Item(29, 4).Value = 1
$ws.Cells.Item(29, 4).NumberFormat = $ws.Cells.Item(2, 4).NumberFormat
$ws.Cells.Item(29, 5).Value = 0.7016
$ws.Cells.Item(29, 6).Value = 0.6174
$ws.Cells.Item(29, 7).Value = "INVPROP"
$ws.Cells.Item(29, 8).Value = "Dhurmeea2016"

# --- Restore the view's active selection to N25 on the pane below the frozen header row ---
$ws.Activate()
$ws.Range("N25").Select()
